# November.xlsx — 11.27 data update + formula/style fixups
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) F131:F139 — convert to one shared formula (F130+D131 pattern)
# ------------------------------------------------------------------
$ws.Range("F131:F139").Formula = "=F130+D131"

# ------------------------------------------------------------------
# 2) New ledger rows 154-164 (data for 11/25-11/27) + row 165's date
# ------------------------------------------------------------------

# -- row 154 --
$ws.Range("B154").Value = 45621
$ws.Range("C154").Value = "午晚饭"
$ws.Range("D154").Value = -15
$ws.Range("E154").Value = "记不太清好像是土豆炸鸡套餐"
$ws.Range("F154").Formula = "=F153+D154"

# -- row 155 --
$ws.Range("B155").Value = 45621
$ws.Range("C155").Value = "入账"
$ws.Range("D155").Value = 10
$ws.Range("E155").Value = "亲爱的勾c发米"
$ws.Range("F155").Formula = "=F154+D155"
$ws.Range("G155").Value = 1000
$ws.Range("H155").Value = 2100

# -- row 156 (小结 / summary, filled style) --
$ws.Range("B156").Value = 45621
$ws.Range("C156").Value = "小结"
$ws.Range("D156").Formula = "=SUM(D150:D155)"
$ws.Range("E156").Value = "*"
$ws.Range("F156").Value = 443.63000000000034
$ws.Range("G156").Value = ""
$ws.Range("H156").Value = ""

# -- row 157 --
$ws.Range("B157").Value = 45622
$ws.Range("C157").Value = "水卡"
$ws.Range("D157").Value = -10
$ws.Range("E157").Value = "洗澡用水卡"
$ws.Range("F157").Formula = "=F156+D157"

# -- row 158 --
$ws.Range("B158").Value = 45622
$ws.Range("C158").Value = "早午饭"
$ws.Range("D158").Value = -12
$ws.Range("E158").Value = "*"
$ws.Range("F158").Formula = "=F157+D158"

# -- row 159 --
$ws.Range("B159").Value = 45622
$ws.Range("C159").Value = "奇妙小蛋挞"
$ws.Range("D159").Value = -5
$ws.Range("E159").Value = "十块钱八个，必须尝尝嘛味"
$ws.Range("F159").Formula = "=F158+D159"

# -- row 160 --
$ws.Range("B160").Value = 45622
$ws.Range("C160").Value = "晚饭"
$ws.Range("D160").Value = -21.78
$ws.Range("E160").Value = "这就是不想出门的下场"
$ws.Range("F160").Formula = "=F159+D160"

# -- row 161 (小结 / summary, filled style) --
$ws.Range("B161").Value = 45622
$ws.Range("C161").Value = "小结"
$ws.Range("D161").Formula = "=SUM(D157:D160)"
$ws.Range("E161").Value = "*"
$ws.Range("F161").Value = 394.85000000000036
$ws.Range("G161").Value = ""
$ws.Range("H161").Value = ""

# -- row 162 --
$ws.Range("B162").Value = 45623
$ws.Range("C162").Value = "无糖小饮料"
$ws.Range("D162").Value = -1
$ws.Range("E162").Value = "中奖瓶盖"
$ws.Range("F162").Formula = "=F161+D162"

# -- row 163 --
$ws.Range("B163").Value = 45623
$ws.Range("C163").Value = "早午饭"
$ws.Range("D163").Value = -7.5
$ws.Range("E163").Value = "9点真是尴尬点，啥都没有，既不承上也不启下"
$ws.Range("F163").Formula = "=F162+D163"

# -- row 164 --
$ws.Range("B164").Value = 45623
$ws.Range("C164").Value = "平帐"
$ws.Range("D164").Value = -5.35
$ws.Range("E164").Value = "使得数据库与实际账面一致"
$ws.Range("F164").Formula = "=F163+D164"

# -- row 165: only the date carries over --
$ws.Range("B165").Value = 45623

# ------------------------------------------------------------------
# 3) Re-apply number/fill formatting that plain .Value writes don't
#    infer on their own, by pasting formats from cells that already
#    carry the right style (reuses existing style indices instead of
#    fabricating new ones).
# ------------------------------------------------------------------

# Date format (style used by existing B-column date cells) for B154:B165
$ws.Range("B153").Copy()
$ws.Range("B154:B165").PasteSpecial(-4122)

# "小结" summary-row fill (style used by row 140) for rows 156 and 161
$ws.Range("C140:H140").Copy()
$ws.Range("C156:H156").PasteSpecial(-4122)
$ws.Range("C161:H161").PasteSpecial(-4122)

# Decimal-fill style used for D160 (matches D153's style)
$ws.Range("D153").Copy()
$ws.Range("D160").PasteSpecial(-4122)
$ws.Range("D160").Value = -21.78

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4) View state: scroll position + active selection
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 141
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D161").Select()

$wb.Save()
